# close #206: Adds support for zero-sum influencing factors
#
# The workbook tracks "influencing factor" proportions per id (rows 3-6)
# across several periods (columns B-H). This change adds support for a
# "DI" (divide by zero / not-applicable) marker: the last data row's
# factors for the "8-2015" and "9-2015" periods are flagged as "DI"
# instead of carrying a stale proportion, and the row above it has its
# equivalent factors zeroed out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-affirm the header labels for the "3-2015"/"4-2015"/"5-2015" period
# columns (B2:D2) as plain period labels.
$ws.Cells.Item(2, 2).Value = "3-2015"
$ws.Cells.Item(2, 3).Value = "4-2015"
$ws.Cells.Item(2, 4).Value = "5-2015"

# Row 5 (id 1100049): zero-out the "8-2015"/"9-2015" (G/H) factors.
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0

# Row 6 (id 1100056): flag the "8-2015"/"9-2015" (G/H) factors as "DI"
# (zero-sum / not computable) instead of a numeric proportion.
$ws.Cells.Item(6, 7).Value = "DI"
$ws.Cells.Item(6, 8).Value = "DI"

# Move the active selection to H6, matching where the edit was made.
$ws.Range("H6").Select()
